$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Remove the five slides (formerly slides 14-18) that separated the
# "Day 1" closing slide from the new "Day 2" closing slide. Delete from the
# highest index down so the remaining indices stay stable while deleting.
# ---------------------------------------------------------------------------
for ($i = 18; $i -ge 14; $i--) {
    $p.Slides.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# The slide that used to be #19 is now #14 (last slide in the deck).
# It carried a small leftover corner-icon picture ("object 3") that the
# author removed, along with renumbering the printed slide-number field.
# ---------------------------------------------------------------------------
$last = $p.Slides.Item($p.Slides.Count)
$last.Shapes.Item("object 3").Delete()
$last.Shapes.Item("Slide Number Placeholder 4").TextFrame.TextRange.Text = "14"

# ---------------------------------------------------------------------------
# Slide 4 had the very same leftover corner-icon picture; remove it too.
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$slide4.Shapes.Item("object 3").Delete()
